# E puts data into long format
# Replace the short species names in column D (Species) with the full
# "Salix <species>" long-format names, for rows 2 through 46 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map short species names -> long format names
$map = @{
    "pulchra"      = "Salix pulchra"
    "richardsonii" = "Salix richardsonii"
    "arctica"      = "Salix arctica"
}

for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Range("D$r")
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# Restore the selection / pane view state recorded in the workbook
$ws.Activate()
$ws.Range("D47").Select()
